# Auto refresh - 16-02-2026 15:26:23.45
#
# The PaceSmart pacing-meta columns (P:S) on "Excel_vs_ML" get reordered
# (DSP_meta, Total_Budget_meta, Flight_Start_Date_meta, Flight_End_Date_meta),
# the per-row meta/prediction columns (P:V) are reset for the new refresh
# cycle, and every campaign's Risk_Level resets to "LOW – Stable" until the
# ML model re-scores them. Feature_Importance and Exec_Summary are
# refreshed with the latest run's numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Excel_vs_ML
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Excel_vs_ML")

# Reorder the meta-column headers in row 1 (P1:S1).
$ws.Range("P1").Value = "DSP_meta"
$ws.Range("Q1").Value = "Total_Budget_meta"
$ws.Range("R1").Value = "Flight_Start_Date_meta"
$ws.Range("S1").Value = "Flight_End_Date_meta"

# Clear the stale per-campaign meta (P:S) and ML prediction (T:V) values
# for every data row - they get repopulated on the next ML scoring pass.
$ws.Range("P2:V36").Clear()

# Every campaign resets to the default "LOW – Stable" risk level for this
# refresh cycle.
for ($r = 2; $r -le 36; $r++) {
    $ws.Range("W$r").Value = "LOW – Stable"
}

# ---------------------------------------------------------------------
# Sheet: Feature_Importance
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Feature_Importance")

$features = @(
    @("Days_Remaining", 0.2787783402760992),
    @("Flight_Days", 0.263697424323086),
    @("Days_Elapsed", 0.1903035724284907),
    @("Total_Budget", 0.1701008344279234),
    @("Gap_Pct", 0.05106095104037819),
    @("Time_Pct", 0.03578353437387945),
    @("DSP_enc", 0.01027534313014306),
    @("Spend_Velocity", 0),
    @("Spend_to_Date", 0),
    @("Acceleration", 0),
    @("Budget_Pct", 0)
)

for ($i = 0; $i -lt $features.Length; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $features[$i][0]
    $ws2.Range("B$row").Value = $features[$i][1]
}

# ---------------------------------------------------------------------
# Sheet: Exec_Summary
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Exec_Summary")
$ws3.Range("B2").Value = 0
$ws3.Range("B3").Value = 0
$ws3.Range("B4").Value = "2026-02-16 09:56 UTC"
